# DRRX Yearly Financials update
# Inserts a new "most recent period" column before column D on the (only)
# worksheet, shifting the existing D:K data one column to the right
# (E:L), and fills the new column D with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# Map of row number -> new value to place in column D (the newly
# inserted "latest period" column). Values were taken from the
# financial statement update.
$newValues = @(
    @{Row=7;   Value=43465},
    @{Row=8;   Value=18600},
    @{Row=9;   Value=4300},
    @{Row=10;  Value=14300},
    @{Row=12;  Value=25500},
    @{Row=13;  Value=0},
    @{Row=14;  Value=0},
    @{Row=15;  Value=0},
    @{Row=17;  Value=42200},
    @{Row=18;  Value=-23600},
    @{Row=20;  Value=900},
    @{Row=21;  Value=-22500},
    @{Row=22;  Value=2600},
    @{Row=23;  Value=-25300},
    @{Row=24;  Value=0},
    @{Row=25;  Value=0},
    @{Row=26;  Value=-25300},
    @{Row=27;  Value=-25300},
    @{Row=28;  Value=0},
    @{Row=29;  Value=0},
    @{Row=30;  Value=0},
    @{Row=31;  Value=0},
    @{Row=32;  Value=-900},
    @{Row=33;  Value=-25300},
    @{Row=34;  Value=0},
    @{Row=35;  Value=-25300},
    @{Row=38;  Value=43465},
    @{Row=41;  Value=31600},
    @{Row=42;  Value=2700},
    @{Row=43;  Value=1800},
    @{Row=44;  Value=3400},
    @{Row=45;  Value=2200},
    @{Row=46;  Value=41700},
    @{Row=47;  Value="NA"},
    @{Row=48;  Value=600},
    @{Row=49;  Value=6400},
    @{Row=50;  Value=0},
    @{Row=51;  Value=0},
    @{Row=52;  Value=1300},
    @{Row=53;  Value=0},
    @{Row=54;  Value=50000},
    @{Row=57;  Value=1600},
    @{Row=58;  Value=0},
    @{Row=59;  Value=6100},
    @{Row=60;  Value=7700},
    @{Row=61;  Value=20500},
    @{Row=62;  Value=1800},
    @{Row=63;  Value=0},
    @{Row=64;  Value=0},
    @{Row=65;  Value=0},
    @{Row=66;  Value=30000},
    @{Row=68;  Value=0},
    @{Row=69;  Value=0},
    @{Row=70;  Value=0},
    @{Row=71;  Value=0},
    @{Row=72;  Value=-468600},
    @{Row=73;  Value=0},
    @{Row=74;  Value=0},
    @{Row=75;  Value=0},
    @{Row=76;  Value=20000},
    @{Row=77;  Value=0},
    @{Row=80;  Value=43465},
    @{Row=81;  Value=-25300},
    @{Row=83;  Value=300},
    @{Row=84;  Value=0},
    @{Row=85;  Value=0},
    @{Row=86;  Value=0},
    @{Row=87;  Value=0},
    @{Row=88;  Value=0},
    @{Row=89;  Value=-19800},
    @{Row=91;  Value=-100},
    @{Row=92;  Value=0},
    @{Row=93;  Value=0},
    @{Row=94;  Value=4500},
    @{Row=96;  Value=0},
    @{Row=97;  Value=0},
    @{Row=98;  Value=0},
    @{Row=99;  Value=0},
    @{Row=100; Value=17500},
    @{Row=101; Value=0},
    @{Row=102; Value=2300}
)

# Rows that have data from column D through column K (now E through L)
# on the sheet. The newly inserted column D cell in each of these rows
# needs to inherit the same number formatting as the rest of the row
# (taken from column E, which now holds what used to be in D).
$dataRows = @(7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,`
              38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,`
              80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102)

foreach ($r in $dataRows) {
    $ws.Range("E$r").Copy() | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

foreach ($item in $newValues) {
    $ws.Range("D$($item.Row)").Value = $item.Value
}
